# Updates the Price (D) and a couple of Volume(1h) (E) cells on Sheet1 to
# match the refreshed crypto-symbol snapshot. The Price column stores its
# numeric-looking values as plain text (inline strings) in the workbook, so
# a naive `.Value = "123.45"` assignment would let Excel auto-coerce the
# cell to a real number (and normalize things like trailing zeros). To keep
# the cell's original text type/style, we briefly force the NumberFormat to
# Text ("@") before writing the value, then restore the original style.

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "243.43"
Set-TextValue $ws.Range("D3") "23.14"
Set-TextValue $ws.Range("D4") "5.419"
Set-TextValue $ws.Range("D5") "0.05941"
Set-TextValue $ws.Range("D6") "3.454"
Set-TextValue $ws.Range("D7") "6.537"
Set-TextValue $ws.Range("D9") "0.9103"
Set-TextValue $ws.Range("D10") "0.1410"
Set-TextValue $ws.Range("D11") "0.07487"
Set-TextValue $ws.Range("D12") "0.03281"
Set-TextValue $ws.Range("D13") "0.03060"
Set-TextValue $ws.Range("D14") "0.09360"
Set-TextValue $ws.Range("D15") "3.857"
Set-TextValue $ws.Range("D16") "0.001558"

Set-TextValue $ws.Range("D18") "0.0005943"
$ws.Range("E18").Value = "17OneONE"

Set-TextValue $ws.Range("D19") "0.006142"
Set-TextValue $ws.Range("D20") "0.004995"
Set-TextValue $ws.Range("D21") "0.0009837"

Set-TextValue $ws.Range("D22") "0.0001101"
$ws.Range("E22").Value = "21NitroExNTXBestin24h"

Set-TextValue $ws.Range("D23") "3.606"

Set-TextValue $ws.Range("D40") "0.03940"
Set-TextValue $ws.Range("D41") "0.006215"
Set-TextValue $ws.Range("D42") "0.1076"
Set-TextValue $ws.Range("D44") "0.007847"
Set-TextValue $ws.Range("D45") "0.00005235"
Set-TextValue $ws.Range("D49") "0.002263"
